# Weekly price update: insert a new "Primera" quality record for Cilantro
# dated 2022-11-29 (serial 44894) ahead of the existing history, shifting
# all subsequent rows down by one (rows 55-123 -> 56-124).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 55 (and everything below it) down one row.
$ws.Rows.Item(55).Insert()

# Populate the newly-inserted row 55 with the new weekly record.
$ws.Range("A55").Value = 7
$ws.Range("B55").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C55").Value = "Ñuble"
$ws.Range("D55").Value = 44894
$ws.Range("E55").Value = 16
$ws.Range("F55").Value = 100112040
$ws.Range("G55").Value = "Cilantro"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 600
$ws.Range("K55").Value = 750
$ws.Range("L55").Value = 850
$ws.Range("M55").Value = 800
$ws.Range("N55").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O55").Value = "Provincia de Diguillín"
$ws.Range("P55").Value = 800
$ws.Range("Q55").Value = 1
$ws.Range("R55").Value = "Hortaliza"
